$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 - this shifts existing rows 14-41 down to 15-42
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with a fresh data record (weekly update)
$ws.Cells.Item(14, 1).Value = 7
$ws.Cells.Item(14, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(14, 3).Value = "Ñuble"
$ws.Cells.Item(14, 4).Value = 44994
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 16
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100101
$ws.Cells.Item(14, 8).Value = "Berries"
$ws.Cells.Item(14, 9).Value = 100101001
$ws.Cells.Item(14, 10).Value = "Arándano (blue)"
$ws.Cells.Item(14, 11).Value = "Sin especificar"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 14).Value = 3000
$ws.Cells.Item(14, 15).Value = 3200
$ws.Cells.Item(14, 16).Value = 3100
$ws.Cells.Item(14, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(14, 19).Value = 1550
$ws.Cells.Item(14, 20).Value = 2
